$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer present in the updated export
# ("Dương Văn Mạnh" and "Lê Văn Mạnh"), which were rows 3 and 4.
# The remaining rows shift up, so:
#   row 2 -> Bùi Thị Như (unchanged person)
#   row 3 -> Nguyễn Hoàng Minh An (was row 5)
#   row 4 -> Đặng Thị Kim Anh (was row 6)
#   row 5 -> Đỗ Thái Dương (was row 7)
#   row 6 -> Đoàn Minh Phương (was row 8)
$ws.Rows("3:4").Delete()

# Every remaining row now shares the corrected "Ngày Nghiệm Thu" date range.
$ws.Range("A2:A6").Value = "5/8/2024"
$ws.Range("B2:B6").Value = "8/12/2024"
$ws.Range("Y2:Y6").Value = "8/12/2024"

# Row 3: Nguyễn Hoàng Minh An - fixed date of birth and recomputed pay figures
$ws.Range("F3").Value = "7/2/1978"
$ws.Range("R3").Value = 43
$ws.Range("S3").Value = 4300000
$ws.Range("T3").Value = "Bốn  Triệu Ba Trăm Không Mươi Nghìn Không Trăm Đồng"
$ws.Range("U3").Value = 430000
$ws.Range("V3").Value = "Bốn Trăm Ba Mươi Nghìn Không Trăm Đồng"
$ws.Range("W3").Value = 3870000
$ws.Range("X3").Value = "Ba  Triệu Tám Trăm Bảy Mươi Nghìn Không Trăm Đồng"

# Row 5: Đỗ Thái Dương - fixed date of birth and recomputed pay figures
$ws.Range("F5").Value = "30/10/2024"
$ws.Range("R5").Value = 65
$ws.Range("S5").Value = 6500000
$ws.Range("T5").Value = "Sáu  Triệu Năm Trăm Không Mươi Nghìn Không Trăm Đồng"
$ws.Range("U5").Value = 650000
$ws.Range("V5").Value = "Sáu Trăm Lăm Nghìn Không Trăm Đồng"
$ws.Range("W5").Value = 5850000
$ws.Range("X5").Value = "Năm  Triệu Tám Trăm Lăm Nghìn Không Trăm Đồng"
